# Update gh-pages to output generated at 456a3b4
# Applies refreshed "F" column (follower/fan count style metric) values
# to the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 3225
$ws1.Range("F5").Value  = 2274
$ws1.Range("F9").Value  = 1041
$ws1.Range("F11").Value = 475
$ws1.Range("F15").Value = 536
$ws1.Range("F16").Value = 8054
$ws1.Range("F17").Value = 351
$ws1.Range("F20").Value = 238
$ws1.Range("F22").Value = 463
$ws1.Range("F23").Value = 547
$ws1.Range("F25").Value = 1139
$ws1.Range("F27").Value = 1862
$ws1.Range("F28").Value = 370
$ws1.Range("F30").Value = 1686
$ws1.Range("F32").Value = 1908
$ws1.Range("F34").Value = 4
$ws1.Range("F37").Value = 287
$ws1.Range("F39").Value = 192
$ws1.Range("F40").Value = 360
$ws1.Range("F42").Value = 230

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 3225
$ws4.Range("F7").Value  = 2274
$ws4.Range("F12").Value = 1041
$ws4.Range("F14").Value = 475
$ws4.Range("F17").Value = 536
$ws4.Range("F18").Value = 8054
$ws4.Range("F19").Value = 351
$ws4.Range("F23").Value = 238
$ws4.Range("F25").Value = 463
$ws4.Range("F26").Value = 547
$ws4.Range("F28").Value = 1139
$ws4.Range("F30").Value = 1862
$ws4.Range("F31").Value = 370
$ws4.Range("F33").Value = 1686
$ws4.Range("F35").Value = 1908
$ws4.Range("F37").Value = 4
$ws4.Range("F40").Value = 287
$ws4.Range("F42").Value = 192
$ws4.Range("F43").Value = 360
$ws4.Range("F49").Value = 230

$wb.Save()
